$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the AR install file references from QA4 to QA2
$ws.Range("A2").Value = "install_zone_fileQA2AR.bat"
$ws.Range("A3").Value = "install_zone_fileQA2AR.bat"
$ws.Range("A4").Value = "install_zone_fileQA2AR.bat"
$ws.Range("A5").Value = "install_zone_fileQA2AR.bat"

# Update the NL install file references from QA4 to QA2
$ws.Range("A6").Value = "install_zone_fileQA2NL.bat"
$ws.Range("A7").Value = "install_zone_fileQA2NL.bat"
$ws.Range("A8").Value = "install_zone_fileQA2NL.bat"
$ws.Range("A9").Value = "install_zone_fileQA2NL.bat"

# Move the active selection to A9 to match final workbook state
$ws.Activate()
$ws.Range("A9").Select()
